# Generate Report for Handback
# Updates the handoff/handback timestamps for the most recently processed
# file (eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md) across the Overview,
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# Row 3 corresponds to eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md
# Column G = "Latest HO Xliff Generate Date" for de-de
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-04 16:53:30"

# --- zh-cn sheet -------------------------------------------------------
# Row 3 corresponds to eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-04 16:53:25"
$wsZhCn.Range("K3").Value = "2016-09-04 16:53:43"

# --- de-de sheet -------------------------------------------------------
# Row 3 corresponds to eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-04 16:53:30"
$wsDeDe.Range("K3").Value = "2016-09-04 16:53:50"
